$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 373.22223
$ws.Range("I38").Value = 294.875
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 884.625
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = -512.625
$ws.Range("N38").Value = -3744
# Row 51
$ws.Range("H51").Value = 6399.8
$ws.Range("I51").Value = 6999.6665
$ws.Range("J51").Value = 5500
$ws.Range("K51").Value = 6999.6665
$ws.Range("L51").Value = 5500
$ws.Range("M51").Value = -6515.6665
$ws.Range("N51").Value = -6468
# Row 58
$ws.Range("H58").Value = 1856.6
$ws.Range("I58").Value = 141.75
$ws.Range("J58").Value = 2999.8333
$ws.Range("K58").Value = 425.25
$ws.Range("L58").Value = 8999.499899999999
$ws.Range("M58").Value = -275.25
$ws.Range("N58").Value = -9299.499899999999
# Row 87
$ws.Range("H87").Value = 99999
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 99999
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 99999
$ws.Range("N87").Value = -102495
# Row 90
$ws.Range("H90").Value = 99999
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 99999
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 299997
$ws.Range("N90").Value = -312477
# Row 116
$ws.Range("H116").Value = 6344.222
$ws.Range("I116").Value = 7366.5
$ws.Range("J116").Value = 4299.6665
$ws.Range("K116").Value = 7366.5
$ws.Range("L116").Value = 4299.6665
$ws.Range("M116").Value = -3924.5
$ws.Range("N116").Value = -11183.6665

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 38
$ws.Range("H38").Value = 9941.25
$ws.Range("I38").Value = 2358.4285
$ws.Range("J38").Value = 63021
$ws.Range("K38").Value = 2358.4285
$ws.Range("L38").Value = 63021
$ws.Range("M38").Value = -1891.4285
$ws.Range("N38").Value = -63955
# Row 39
$ws.Range("H39").Value = 3237.1428
$ws.Range("I39").Value = 2943.3333
$ws.Range("J39").Value = 5000
$ws.Range("K39").Value = 2943.3333
$ws.Range("L39").Value = 5000
$ws.Range("M39").Value = -2423.3333
$ws.Range("N39").Value = -6040
# Row 45
$ws.Range("H45").Value = 3304.75
$ws.Range("I45").Value = 1785.5
$ws.Range("J45").Value = 4824
$ws.Range("K45").Value = 1785.5
$ws.Range("L45").Value = 4824
$ws.Range("M45").Value = -1408.5
$ws.Range("N45").Value = -5578
# Row 47
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").ClearContents()
$ws.Range("N47").Value = 0
# Row 50
$ws.Range("H50").Value = 1334.3334
$ws.Range("I50").Value = 450
$ws.Range("J50").Value = 1776.5
$ws.Range("K50").Value = 450
$ws.Range("L50").Value = 1776.5
$ws.Range("M50").Value = 264
$ws.Range("N50").Value = -3204.5
# Row 53
$ws.Range("H53").Value = 199997
$ws.Range("I53").Value = 199997
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 199997
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -199315

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 59.25
$ws.Range("I7").Value = 59.11111
$ws.Range("J7").Value = 59.666668
$ws.Range("K7").Value = 59.11111
$ws.Range("L7").Value = 59.666668
$ws.Range("M7").Value = 53.88889
$ws.Range("N7").Value = -285.666668
# Row 58
$ws.Range("H58").Value = 4727.222
$ws.Range("I58").Value = 1478
$ws.Range("J58").Value = 7326.6
$ws.Range("K58").Value = 1478
$ws.Range("L58").Value = 7326.6
$ws.Range("M58").Value = -1275
$ws.Range("N58").Value = -7732.6
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = 0
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = 0
# Row 107
$ws.Range("H107").Value = 479.85715
$ws.Range("I107").Value = 552
$ws.Range("J107").Value = 299.5
$ws.Range("K107").Value = 552
$ws.Range("L107").Value = 299.5
$ws.Range("M107").Value = 1368
$ws.Range("N107").Value = -4139.5
# Row 136
$ws.Range("H136").Value = 4727.222
$ws.Range("I136").Value = 1478
$ws.Range("J136").Value = 7326.6
$ws.Range("K136").Value = 4434
$ws.Range("L136").Value = 21979.8
$ws.Range("M136").Value = -1884
$ws.Range("N136").Value = -27079.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 39.875
$ws.Range("I7").Value = 35.666668
$ws.Range("J7").Value = 52.5
$ws.Range("K7").Value = 107.000004
$ws.Range("L7").Value = 157.5
$ws.Range("M7").Value = 4.999995999999996
$ws.Range("N7").Value = -381.5
# Row 25
$ws.Range("H25").Value = 37
$ws.Range("I25").Value = 6
$ws.Range("J25").Value = 52.5
$ws.Range("K25").Value = 18
$ws.Range("L25").Value = 157.5
$ws.Range("M25").Value = 151
$ws.Range("N25").Value = -495.5
# Row 30
$ws.Range("H30").Value = 37
$ws.Range("I30").Value = 6
$ws.Range("J30").Value = 52.5
$ws.Range("K30").Value = 18
$ws.Range("L30").Value = 157.5
$ws.Range("M30").Value = 84
$ws.Range("N30").Value = -361.5
# Row 34
$ws.Range("H34").Value = 3149.6
$ws.Range("I34").Value = 500
$ws.Range("J34").Value = 3812
$ws.Range("K34").Value = 1500
$ws.Range("L34").Value = 11436
$ws.Range("M34").Value = -1416
$ws.Range("N34").Value = -11604
# Row 39
$ws.Range("H39").Value = 7200
$ws.Range("I39").Value = 800
$ws.Range("J39").Value = 9333.333000000001
$ws.Range("K39").Value = 2400
$ws.Range("L39").Value = 27999.999
$ws.Range("M39").Value = -2106
$ws.Range("N39").Value = -28587.999
# Row 55
$ws.Range("H55").Value = 6160.2
$ws.Range("I55").Value = 1552
$ws.Range("J55").Value = 7312.25
$ws.Range("K55").Value = 4656
$ws.Range("L55").Value = 21936.75
$ws.Range("M55").Value = -4479
$ws.Range("N55").Value = -22290.75
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("N62").Value = 0
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("N65").Value = 0

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 40.375
$ws.Range("I2").Value = 23.666666
$ws.Range("J2").Value = 50.4
$ws.Range("K2").Value = 23.666666
$ws.Range("L2").Value = 50.4
$ws.Range("M2").Value = 89.33333400000001
$ws.Range("N2").Value = -276.4
# Row 36
$ws.Range("H36").Value = 5450
$ws.Range("I36").Value = 900
$ws.Range("J36").Value = 10000
$ws.Range("K36").Value = 900
$ws.Range("L36").Value = 10000
$ws.Range("M36").Value = -415
$ws.Range("N36").Value = -10970
# Row 102
$ws.Range("H102").Value = 1381.9231
$ws.Range("I102").Value = 1578.8182
$ws.Range("J102").Value = 299
$ws.Range("K102").Value = 1578.8182
$ws.Range("L102").Value = 299
$ws.Range("M102").Value = 43.18180000000007
$ws.Range("N102").Value = -3543
# Row 126
$ws.Range("H126").Value = 6387.25
$ws.Range("I126").Value = 5616.8
$ws.Range("J126").Value = 7671.3335
$ws.Range("K126").Value = 16850.4
$ws.Range("L126").Value = 23014.0005
$ws.Range("M126").Value = -14380.4
$ws.Range("N126").Value = -27954.0005

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4764.8
$ws.Range("I7").Value = 3135.125
$ws.Range("J7").Value = 6627.2856
$ws.Range("K7").Value = 3135.125
$ws.Range("L7").Value = 6627.2856
$ws.Range("M7").Value = -3023.125
$ws.Range("N7").Value = -6851.2856
# Row 22
$ws.Range("H22").Value = 857.5
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 936.25
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 936.25
$ws.Range("M22").Value = -405
$ws.Range("N22").Value = -1526.25
# Row 27
$ws.Range("H27").Value = 857.5
$ws.Range("I27").Value = 700
$ws.Range("J27").Value = 936.25
$ws.Range("K27").Value = 700
$ws.Range("L27").Value = 936.25
$ws.Range("M27").Value = -593
$ws.Range("N27").Value = -1150.25
# Row 46
$ws.Range("H46").Value = 5306.35
$ws.Range("I46").Value = 3569.6667
$ws.Range("J46").Value = 6727.273
$ws.Range("K46").Value = 3569.6667
$ws.Range("L46").Value = 6727.273
$ws.Range("M46").Value = -3381.6667
$ws.Range("N46").Value = -7103.273
# Row 55
$ws.Range("H55").Value = 1241.1111
$ws.Range("I55").Value = 1352.5
$ws.Range("J55").Value = 350
$ws.Range("K55").Value = 1352.5
$ws.Range("L55").Value = 350
$ws.Range("M55").Value = -1179.5
$ws.Range("N55").Value = -696
# Row 126
$ws.Range("H126").Value = 4764.8
$ws.Range("I126").Value = 3135.125
$ws.Range("J126").Value = 6627.2856
$ws.Range("K126").Value = 9405.375
$ws.Range("L126").Value = 19881.8568
$ws.Range("M126").Value = -6935.375
$ws.Range("N126").Value = -24821.8568
# Row 136
$ws.Range("H136").Value = 2877.7778
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 2987.5
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 8962.5
$ws.Range("M136").Value = -3450
$ws.Range("N136").Value = -14062.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 45125.668
$ws.Range("I41").Value = 40000
$ws.Range("J41").Value = 47688.5
$ws.Range("K41").Value = 40000
$ws.Range("L41").Value = 47688.5
$ws.Range("M41").Value = -39610
$ws.Range("N41").Value = -48468.5
# Row 126
$ws.Range("H126").Value = 3768.3044
$ws.Range("I126").Value = 1641.5
$ws.Range("J126").Value = 6088.4546
$ws.Range("K126").Value = 4924.5
$ws.Range("L126").Value = 18265.3638
$ws.Range("M126").Value = -2454.5
$ws.Range("N126").Value = -23205.3638
# Row 132
$ws.Range("H132").Value = 3935.5
$ws.Range("I132").Value = 2786.625
$ws.Range("J132").Value = 5467.3335
$ws.Range("K132").Value = 8359.875
$ws.Range("L132").Value = 16402.0005
$ws.Range("M132").Value = -5829.875
$ws.Range("N132").Value = -21462.0005
